# The document currently ends with an empty paragraph followed by the
# sectPr. The diff keeps that trailing empty paragraph untouched and then
# appends a run of new paragraphs (starting with another empty one, then the
# new "How to edit files in GCAF_git..." section) right before the sectPr.
#
# Word's Range.InsertXML() splices raw WordprocessingML at a collapsed
# range, but when that range sits at the very end of the document it fuses
# the first inserted <w:p> into the paragraph that currently owns the end
# mark (dropping that paragraph's own rsid attributes). Calling
# InsertParagraphAfter() first manufactures a brand-new, untouched empty
# paragraph at the end of the story, so the existing last paragraph is left
# exactly as-is and our XML block lands entirely after it.

$d = $word.ActiveDocument

$endOfDoc = $d.Paragraphs.Last.Range
$endOfDoc.Collapse(0)
[void]$endOfDoc.InsertParagraphAfter()

$insertionRange = $d.Paragraphs.Last.Range
$insertionRange.Collapse(0)

$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">How to edit files in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>GCAF_git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> and update (push) them to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Edit the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pkg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> files, and re-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tarball</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> them to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tar.gz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file with the exact name GrowthCurveAnalysis_0.0.2.tar.gz (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>REtarballing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is important, so when folks install in R, they’re getting your latest edits!)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">In terminal: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>cd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ~/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GCAF_git</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GrowthCurveAnalysis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> GrowthCurveAnalysis_0.0.2.tar.gz </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> commit –m ‘Note about your addition/edit’</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> push </w:t></w:r></w:p>
'@

[void]$insertionRange.InsertXML($newContentXml)
